$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")
$summary = $wb.Worksheets.Item("Summary")

# New data rows to append below the existing row 11 (A11 already holds "3").
# Each entry: row number -> column letter -> text value.
$rowsData = @(
    @{ Row = 11; A = $null; C = "556_马尾松_undefined_undefined_1bunch"; F = "30" },
    @{ Row = 12; A = $null; C = "345_天竺少女_Cryptomeria`nKashiwaba_undefined_1bunch"; F = "10" },
    @{ Row = 13; A = "4";   C = "342_南天竹红_undefined_Nandina domestica Thunb._1bunch"; F = "10" },
    @{ Row = 14; A = $null; C = "343_南天竹带果_undefined_Nandina domestica Thunb._1bunch"; F = "10" },
    @{ Row = 15; A = $null; C = "321_雪柳叶_Spiraea  leaves_undefined_1bunch"; F = "30" },
    @{ Row = 16; A = $null; C = "106_绣球单瓣粉_Hydrangea Pink S_Hydrangea L._1stem"; F = "50" },
    @{ Row = 17; A = $null; C = "110_绣球单瓣浅蓝_Hydrangea Light Blue S_Hydrangea L._1stem"; F = "40" },
    @{ Row = 18; A = $null; C = "11_香槟洋桔梗_Champagne Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"; F = "10" },
    @{ Row = 19; A = $null; C = "624_多丁白_undefined_undefined_1bunch"; F = "5" },
    @{ Row = 20; A = "5";   C = "475_诺贝松_undefined_undefined_1bunch"; F = "20" },
    @{ Row = 21; A = $null; C = "597_尤加利叶小叶_undefined_undefined_1bunch"; F = $null }
)

# First, make sure every cell that will receive a numeric-looking text value
# is formatted as Text ("@") BEFORE the value is assigned, otherwise Excel
# auto-converts the numeric-looking string into a real number.
foreach ($entry in $rowsData) {
    $r = $entry.Row
    if ($entry.A) { $ws.Range("A$r").NumberFormat = "@" }
    if ($entry.F) { $ws.Range("F$r").NumberFormat = "@" }
}

foreach ($entry in $rowsData) {
    $r = $entry.Row
    if ($entry.A) { $ws.Range("A$r").Value = $entry.A }
    if ($entry.C) { $ws.Range("C$r").Value = $entry.C }
    if ($entry.F) { $ws.Range("F$r").Value = $entry.F }
}

# Row 12 contains an embedded newline; Excel auto-expands the row height
# when that happens. Restore the row to its default (auto) height so the
# saved file doesn't carry a stray custom row height.
$ws.Rows.Item(12).AutoFit()

# Update the Summary sheet's concatenated "Number" digest string (column G)
# to include the newly added rows' quantities.
$summary.Range("G2").NumberFormat = "@"
$summary.Range("G2").Value = "016111210101614101030101010305040105200"
